$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row, Coin (B), Link (C), Price (D), Volume(1h) (E)
$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "60.488.55", "  -5.06%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.440.47", "  -6.68%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.27%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "540.68", "  -6.00%  "),
    @(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "143.27", "  -8.41%  "),
    @(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.18%  "),
    @(8, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.603", "  -3.22%  "),
    @(9, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "2.439.19", "  -6.64%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.105", "  -11.46%  "),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.153", "  -2.04%  "),
    @(12, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "5.31", "  -9.21%  "),
    @(13, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.348", "  -8.53%  "),
    @(14, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "25.61", "  -9.45%  "),
    @(15, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.895.69", "  -6.01%  "),
    @(16, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "60.472.44", "  -4.90%  "),
    @(17, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000160", "  -10.67%  "),
    @(18, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.467.22", "  -5.63%  "),
    @(19, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "10.96", "  -8.93%  "),
    @(20, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.87", "  -9.78%  "),
    @(21, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.12", "  -9.27%  "),
    @(22, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "315.16", "  -8.29%  "),
    @(23, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.05%  "),
    @(24, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "62.79", "  -7.07%  "),
    @(25, "SuiNetwork", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui", "1.71", "  -5.35%  "),
    @(26, "WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "2.606.29", "  -4.43%  "),
    @(27, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "1.00", "  +0.20%  "),
    @(28, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0946", "  -13.30%  "),
    @(29, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.46", "  -7.81%  "),
    @(30, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "8.20", "  -10.64%  "),
    @(31, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "523.14", "  -10.95%  "),
    @(32, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.49", "  -5.35%  "),
    @(33, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.146", "  -9.85%  "),
    @(34, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.87", "  -9.19%  "),
    @(35, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.54", "  -11.42%  "),
    @(36, "RenderToken", "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render", "5.71", "  -13.47%  "),
    @(37, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.00", "  +0.28%  "),
    @(38, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "4.77", "  -11.15%  "),
    @(39, "PolygonEcosystemToken", "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol", "0.371", "  -8.02%  "),
    @(40, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.17", "  -8.14%  "),
    @(41, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "143.11", "  -7.31%  "),
    @(42, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "0.999", "  +0.01%  "),
    @(43, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.67", "  -11.03%  "),
    @(44, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "39.93", "  -3.61%  "),
    @(45, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "2.24", "  -11.59%  "),
    @(46, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "145.12", "  -7.85%  "),
    @(47, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.53", "  -9.65%  "),
    @(48, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "20.58", "  -13.43%  "),
    @(49, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0523", "  -11.43%  "),
    @(50, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.0931", "  -7.16%  "),
    @(51, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.575", "  -8.83%  "),
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    # Column D holds price strings (e.g. "1.00", "0.0000160") that must stay text,
    # so force text format before assignment to avoid Excel auto-converting to a number,
    # then restore the default "Normal" style so no extra formatting is introduced.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item[3]
    $dCell.Style = "Normal"
    $ws.Cells.Item($r, 5).Value = $item[4]
}

Write-Host "Updated $($data.Count) rows"
